$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PACH")

Write-Output "BEFORE MOVE:"
Write-Output $ws.Range("A25").Value()
Write-Output $ws.Range("A65").Value()

$src = $ws.Rows("65:79")
$src.Cut($ws.Rows("25:25")) | Out-Null

Write-Output "AFTER MOVE:"
Write-Output $ws.Range("A25").Value()
Write-Output $ws.Range("A40").Value()
Write-Output $ws.Range("A79").Value()
